$wb = $excel.ActiveWorkbook

# The "Trust Vs Violations" sheet is the second sheet in the workbook.
$ws = $wb.Worksheets.Item("Trust Vs Violations")

# --- Add the new Log-violations (F) and duplicated trustcompanies (G) columns ---
$ws.Range("F1").Value = "Log violations "
$ws.Range("G1").Value = "trustcompanies"

for ($r = 2; $r -le 68; $r++) {
    $ws.Range("F$r").Formula = "=LOG(C$r +0.1)"
    $ws.Range("G$r").Value = $ws.Range("D$r").Value2
}

# --- Make "Trust Vs Violations" the active/selected sheet & cell ---
$ws.Activate() | Out-Null
$ws.Range("E11").Select() | Out-Null

# --- Move / resize the scatter-plot chart on this sheet ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 435.3125
$co.Top = 9.25
$co.Width = 638.4374606299214
$co.Height = 234.0
